$wb = $excel.ActiveWorkbook

# --- Update the "Login" sheet test data -----------------------------------
# Three separate rows ("ScreeningHappyPath_ErWork", "ScreeningHappyPath_Sup",
# "ScreeningHappyPath_Staff") are consolidated into three iterations of the
# same "ScreeningHappyPath" test case (Iteration 1, 2, 3).
$wsLogin = $wb.Worksheets.Item("Login")

$wsLogin.Range("A2").Value = "ScreeningHappyPath"
$wsLogin.Range("A3").Value = "ScreeningHappyPath"
$wsLogin.Range("A4").Value = "ScreeningHappyPath"

# B2 already holds Iteration "1" and stays that way.
# B3 / B4 become text "2" / "3" (typed with a leading apostrophe so Excel
# keeps them as text, matching the quote-prefixed style used in the sheet).
$wsLogin.Range("B3").Value = "'2"
$wsLogin.Range("B4").Value = "'3"

# --- Update sheet selection / active sheet ---------------------------------
# Previously "InitialScreening" was the active tab with G5 selected; now
# "Login" is the active tab (with C8 selected) and "InitialScreening" keeps
# B2 selected for when it's revisited.
$wsInitialScreening = $wb.Worksheets.Item("InitialScreening")
$wsInitialScreening.Range("B2").Select() | Out-Null

$wsLogin.Select() | Out-Null
$wsLogin.Range("C8").Select() | Out-Null
